$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26 ---
$ws.Range("A26").Value = 42580
$ws.Range("A26").NumberFormat = $ws.Range("A18").NumberFormat

$ws.Range("B26").Value = "8-core Macpro"
$ws.Range("B26").NumberFormat = $ws.Range("B18").NumberFormat

$ws.Range("C26").Value = "Sim"

$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 1200
$ws.Range("F26").Value = 10000

$ws.Range("H26").Value = 4.3
$ws.Range("H26").HorizontalAlignment = -4152

$ws.Range("I26").Value = 5

# --- Row 27 ---
$ws.Range("A27").Value = 42580
$ws.Range("A27").NumberFormat = $ws.Range("A18").NumberFormat

$ws.Range("B27").Value = "8-core Macpro"
$ws.Range("B27").NumberFormat = $ws.Range("B18").NumberFormat

$ws.Range("C27").Value = "Sim"

$ws.Range("D27").Value = 28
$ws.Range("E27").Value = 1200
$ws.Range("F27").Value = 10000

$ws.Range("H27").Value = 5.0999999999999996
$ws.Range("H27").HorizontalAlignment = -4152

$ws.Range("I27").Style = "Normal"
$ws.Range("I27").Value = 5

$ws.Range("J27").Value = 1
$ws.Range("J27").ClearContents()
$ws.Range("J27").Style = "Normal"

$ws.Range("M27").Value = 1
$ws.Range("M27").ClearContents()
$ws.Range("M27").Style = "Normal"

# --- Row 28 ---
$ws.Range("A28").Value = 42580
$ws.Range("A28").NumberFormat = $ws.Range("A18").NumberFormat

$ws.Range("B28").Value = "8-core Macpro"
$ws.Range("B28").NumberFormat = $ws.Range("B18").NumberFormat

$ws.Range("C28").Value = "Sim"

$ws.Range("D28").Value = 29
$ws.Range("E28").Value = 1200
$ws.Range("F28").Value = 10000

$ws.Range("H28").Value = 25.2
$ws.Range("H28").HorizontalAlignment = -4152

$ws.Range("I28").Style = "Normal"
$ws.Range("I28").Value = 5

$ws.Range("J28").Value = 1
$ws.Range("J28").ClearContents()
$ws.Range("J28").Style = "Normal"

$ws.Range("M28").Value = 1
$ws.Range("M28").ClearContents()
$ws.Range("M28").Style = "Normal"

# --- Row 29 ---
$ws.Range("A29").Value = 42580
$ws.Range("A29").NumberFormat = $ws.Range("A18").NumberFormat

$ws.Range("B29").Value = "8-core Macpro"
$ws.Range("B29").NumberFormat = $ws.Range("B18").NumberFormat

$ws.Range("C29").Value = "Sim"

$ws.Range("D29").Value = 31
$ws.Range("E29").Value = 1200
$ws.Range("F29").Value = 10000

$ws.Range("H29").Value = 25.2
$ws.Range("H29").HorizontalAlignment = -4152

$ws.Range("I29").Style = "Normal"
$ws.Range("I29").Value = 5

$ws.Range("J29").Value = 1
$ws.Range("J29").ClearContents()
$ws.Range("J29").Style = "Normal"

$ws.Range("M29").Value = 1
$ws.Range("M29").ClearContents()
$ws.Range("M29").Style = "Normal"

# --- Shared formula G26:G29 = E*F ---
$ws.Range("G26:G29").Formula = "=E26*F26"

# --- Update selected cell in sheet view ---
$ws.Range("I24").Select() | Out-Null
